# Weekly driver report update for 2025-04-20
# Updates the "Good Drivers" table (rows 12-17) on the active sheet:
#  - re-sorted order of adapter/driver rows
#  - refreshed Client/Total sample counts
#  - refreshed Good Roaming Calculation (%) values
#  - refreshed Driver Vintage dates (stored as literal text, matching the
#    source report's plain-text date column)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the "Driver Vintage" column (E) as literal text so dates like
# "2024-11-10" aren't reinterpreted as date serials by Excel's
# type-inference on assignment.
$ws.Range("E12:E17").NumberFormat = "@"

# Row 12
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B12").Value = 445055
$ws.Range("D12").Value = 99.90000000000001
$ws.Range("E12").Value = "2024-11-10"

# Row 13
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B13").Value = 77849
$ws.Range("D13").Value = 99.90000000000001
$ws.Range("E13").Value = "2021-08-18"

# Row 14
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B14").Value = 34244
$ws.Range("D14").Value = 100
$ws.Range("E14").Value = "2021-04-27"

# Row 15
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B15").Value = 59673
$ws.Range("D15").Value = 100
$ws.Range("E15").Value = "2020-08-05"

# Row 16
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B16").Value = 113652
$ws.Range("D16").Value = 100
$ws.Range("E16").Value = "2020-01-06"

# Row 17
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B17").Value = 56018
$ws.Range("D17").Value = 100
$ws.Range("E17").Value = "2019-12-14"
